$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 1511.0555
$ws.Range("I39").Value = 989
$ws.Range("J39").Value = 2555.1667
$ws.Range("K39").Value = 2967
$ws.Range("L39").Value = 7665.500100000001
$ws.Range("M39").Value = -2671
$ws.Range("N39").Value = -8257.500100000001
$ws.Range("H51").Value = 8252.111000000001
$ws.Range("J51").Value = 8938.429
$ws.Range("L51").Value = 8938.429
$ws.Range("N51").Value = -9906.429
$ws.Range("H88").Value = 2409
$ws.Range("I88").Value = 834
$ws.Range("J88").Value = 2704.3125
$ws.Range("K88").Value = 834
$ws.Range("L88").Value = 2704.3125
$ws.Range("M88").Value = -428
$ws.Range("N88").Value = -3516.3125
$ws.Range("H91").Value = 2409
$ws.Range("I91").Value = 834
$ws.Range("J91").Value = 2704.3125
$ws.Range("K91").Value = 834
$ws.Range("L91").Value = 2704.3125
$ws.Range("M91").Value = 570
$ws.Range("N91").Value = -5512.3125
$ws.Range("H107").Value = 1305.6875
$ws.Range("I107").Value = 1380.2307
$ws.Range("K107").Value = 1380.2307
$ws.Range("M107").Value = 539.7692999999999
$ws.Range("H132").Value = 7163.8726
$ws.Range("I132").Value = 5324.625
$ws.Range("K132").Value = 15973.875
$ws.Range("M132").Value = -13443.875
$ws.Range("H138").Value = 7454.836
$ws.Range("I138").Value = 10128.637
$ws.Range("J138").Value = 6866.6
$ws.Range("K138").Value = 30385.911
$ws.Range("L138").Value = 20599.8
$ws.Range("M138").Value = -25245.911
$ws.Range("N138").Value = -30879.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 497.25
$ws.Range("I4").Value = 490
$ws.Range("K4").Value = 490
$ws.Range("M4").Value = -374
$ws.Range("H32").Value = 31221.182
$ws.Range("I32").Value = 24776.857
$ws.Range("K32").Value = 24776.857
$ws.Range("M32").Value = -24489.857
$ws.Range("H61").Value = 7352.7334
$ws.Range("I61").Value = 6968.9165
$ws.Range("K61").Value = 6968.9165
$ws.Range("M61").Value = -6756.9165
$ws.Range("H110").Value = 4651.25
$ws.Range("J110").Value = 3900
$ws.Range("L110").Value = 3900
$ws.Range("N110").Value = -7990
$ws.Range("H136").Value = 7352.7334
$ws.Range("I136").Value = 6968.9165
$ws.Range("K136").Value = 20906.7495
$ws.Range("M136").Value = -18356.7495

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H21").Value = 9704.200000000001
$ws.Range("J21").Value = 9704.200000000001
$ws.Range("L21").Value = 9704.200000000001
$ws.Range("N21").Value = -10176.2
$ws.Range("H82").Value = 37372.5
$ws.Range("I82").Value = 25198.2
$ws.Range("J82").Value = 57663
$ws.Range("K82").Value = 25198.2
$ws.Range("L82").Value = 57663
$ws.Range("M82").Value = -24815.2
$ws.Range("N82").Value = -58429
$ws.Range("H85").Value = 37372.5
$ws.Range("I85").Value = 25198.2
$ws.Range("J85").Value = 57663
$ws.Range("K85").Value = 25198.2
$ws.Range("L85").Value = 57663
$ws.Range("M85").Value = -23872.2
$ws.Range("N85").Value = -60315

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 122.73333
$ws.Range("J7").Value = 126.57143
$ws.Range("L7").Value = 126.57143
$ws.Range("N7").Value = -352.57143
$ws.Range("H58").Value = 62273.766
$ws.Range("J58").Value = 3343
$ws.Range("L58").Value = 3343
$ws.Range("N58").Value = -3749
$ws.Range("H86").Value = 121980.5
$ws.Range("I86").Value = 203963
$ws.Range("J86").Value = 39998
$ws.Range("K86").Value = 203963
$ws.Range("L86").Value = 39998
$ws.Range("M86").Value = -202840
$ws.Range("N86").Value = -42244
$ws.Range("H89").Value = 121980.5
$ws.Range("I89").Value = 203963
$ws.Range("J89").Value = 39998
$ws.Range("K89").Value = 1019815
$ws.Range("L89").Value = 199990
$ws.Range("M89").Value = -1014199
$ws.Range("N89").Value = -211222
$ws.Range("H134").Value = 32594.943
$ws.Range("I134").Value = 37084.234
$ws.Range("J134").Value = 5659.2
$ws.Range("K134").Value = 111252.702
$ws.Range("L134").Value = 16977.6
$ws.Range("M134").Value = -108717.702
$ws.Range("N134").Value = -22047.6
$ws.Range("H136").Value = 62273.766
$ws.Range("J136").Value = 3343
$ws.Range("L136").Value = 10029
$ws.Range("N136").Value = -15129
$ws.Range("H139").Value = 79999
$ws.Range("J139").Value = 79999
$ws.Range("L139").Value = 79999
$ws.Range("N139").Value = -90279

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 5293.778
$ws.Range("I141").Value = 5293.778
$ws.Range("K141").Value = 15881.334
$ws.Range("M141").Value = -10701.334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4082.6667
$ws.Range("I70").Value = 3599.5
$ws.Range("J70").Value = 4324.25
$ws.Range("K70").Value = 3599.5
$ws.Range("L70").Value = 4324.25
$ws.Range("M70").Value = -3329.5
$ws.Range("N70").Value = -4864.25
$ws.Range("H73").Value = 4082.6667
$ws.Range("I73").Value = 3599.5
$ws.Range("J73").Value = 4324.25
$ws.Range("K73").Value = 3599.5
$ws.Range("L73").Value = 4324.25
$ws.Range("M73").Value = -2663.5
$ws.Range("N73").Value = -6196.25
$ws.Range("H80").Value = 3999
$ws.Range("I80").Value = 3498.5
$ws.Range("K80").Value = 3498.5
$ws.Range("M80").Value = -2500.5
$ws.Range("H83").Value = 3999
$ws.Range("I83").Value = 3498.5
$ws.Range("K83").Value = 17492.5
$ws.Range("M83").Value = -12500.5
$ws.Range("H132").Value = 193625.73
$ws.Range("I132").Value = 125164.89
$ws.Range("J132").Value = 501699.5
$ws.Range("K132").Value = 375494.67
$ws.Range("L132").Value = 1505098.5
$ws.Range("M132").Value = -372964.67
$ws.Range("N132").Value = -1510158.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 729467.7
$ws.Range("I7").Value = 1451528.4
$ws.Range("J7").Value = 7407
$ws.Range("K7").Value = 1451528.4
$ws.Range("L7").Value = 7407
$ws.Range("M7").Value = -1451416.4
$ws.Range("N7").Value = -7631
$ws.Range("H40").Value = 4026.2354
$ws.Range("I40").Value = 3204.3333
$ws.Range("J40").Value = 5998.8
$ws.Range("K40").Value = 3204.3333
$ws.Range("L40").Value = 5998.8
$ws.Range("M40").Value = -3068.3333
$ws.Range("N40").Value = -6270.8
$ws.Range("H68").Value = 5479.4614
$ws.Range("I68").Value = 4318
$ws.Range("J68").Value = 6205.375
$ws.Range("K68").Value = 4318
$ws.Range("L68").Value = 6205.375
$ws.Range("M68").Value = -3569
$ws.Range("N68").Value = -7703.375
$ws.Range("H71").Value = 5479.4614
$ws.Range("I71").Value = 4318
$ws.Range("J71").Value = 6205.375
$ws.Range("K71").Value = 21590
$ws.Range("L71").Value = 31026.875
$ws.Range("M71").Value = -17846
$ws.Range("N71").Value = -38514.875
$ws.Range("H126").Value = 729467.7
$ws.Range("I126").Value = 1451528.4
$ws.Range("J126").Value = 7407
$ws.Range("K126").Value = 4354585.199999999
$ws.Range("L126").Value = 22221
$ws.Range("M126").Value = -4352115.199999999
$ws.Range("N126").Value = -27161

